$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Force the A2:B7 data block (order/delivery date strings) to Text format so
# the "YYYY-MM-DD" literals are kept as plain text instead of being auto-parsed
# into date serials when assigned through .Value (matches source template,
# which stores them as plain strings, not dates). Scoped to just the rows we
# touch so no other cell/style is disturbed.
$ws1.Range("A2:B7").NumberFormat = "@"

# --- Row 1 headers: rewrite A1:P1, strip the bold/bordered header style ---
$headers = @(
  "발주일자", "납기일자", "거래처명", "거래처 이메일", "납품처명", "납품처 이메일", "프로젝트명", "대분류", "중분류", "소분류", "품목명", "규격", "수량", "단가", "총금액", "비고"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $cell = $ws1.Cells.Item(1, $i + 1)
  $cell.Value = $headers[$i]
  $cell.Style = "Normal"
}

# --- Row 2 ---
$ws1.Cells.Item(2, 1).Value = "2025-08-22"
$ws1.Cells.Item(2, 2).Value = "2025-08-29"
$ws1.Cells.Item(2, 3).Value = "티에스이앤씨"
$ws1.Cells.Item(2, 4).Value = "티에스이앤씨@example.com"
$ws1.Cells.Item(2, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(2, 6).Value = "delivery@example.com"
$ws1.Cells.Item(2, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(2, 8).Value = "6. 안전관리비"
$ws1.Cells.Item(2, 9).Value = "1) 안전장비"
$ws1.Cells.Item(2, 10).Value = "기타"
$ws1.Cells.Item(2, 11).Value = "안전 1차 - 안전벨트"
$ws1.Cells.Item(2, 12).Value = "KS규격-1"
$ws1.Cells.Item(2, 13).Value = 5
$ws1.Cells.Item(2, 14).Value = 37000
$ws1.Cells.Item(2, 15).Value = 203500
$ws1.Cells.Item(2, 16).ClearContents()

# --- Row 3 ---
$ws1.Cells.Item(3, 1).Value = "2025-09-17"
$ws1.Cells.Item(3, 2).Value = "2025-10-03"
$ws1.Cells.Item(3, 3).Value = "티에스이앤씨"
$ws1.Cells.Item(3, 4).Value = "티에스이앤씨@example.com"
$ws1.Cells.Item(3, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(3, 6).Value = "delivery@example.com"
$ws1.Cells.Item(3, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(3, 8).Value = "6. 안전관리비"
$ws1.Cells.Item(3, 9).Value = "1) 안전장비"
$ws1.Cells.Item(3, 10).Value = "기타"
$ws1.Cells.Item(3, 11).Value = "화재감시자 안전모"
$ws1.Cells.Item(3, 12).Value = "KS규격-2"
$ws1.Cells.Item(3, 13).Value = 1
$ws1.Cells.Item(3, 14).Value = 6000
$ws1.Cells.Item(3, 15).Value = 6600
$ws1.Cells.Item(3, 16).ClearContents()

# --- Row 4 ---
$ws1.Cells.Item(4, 1).Value = "2025-08-29"
$ws1.Cells.Item(4, 2).Value = "2025-10-01"
$ws1.Cells.Item(4, 3).Value = "티에스이앤씨"
$ws1.Cells.Item(4, 4).Value = "티에스이앤씨@example.com"
$ws1.Cells.Item(4, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(4, 6).Value = "delivery@example.com"
$ws1.Cells.Item(4, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(4, 8).Value = "2. 부자재비"
$ws1.Cells.Item(4, 9).Value = "2) 창호"
$ws1.Cells.Item(4, 10).Value = "기타"
$ws1.Cells.Item(4, 11).Value = "3차 - 스크류 (둥근머리 8*25)"
$ws1.Cells.Item(4, 12).Value = "KS규격-3"
$ws1.Cells.Item(4, 13).Value = 500
$ws1.Cells.Item(4, 14).Value = 19
$ws1.Cells.Item(4, 15).Value = 10450
$ws1.Cells.Item(4, 16).ClearContents()

# --- Row 5 ---
$ws1.Cells.Item(5, 1).Value = "2025-08-29"
$ws1.Cells.Item(5, 2).Value = "2025-09-22"
$ws1.Cells.Item(5, 3).Value = "티에스이앤씨"
$ws1.Cells.Item(5, 4).Value = "티에스이앤씨@example.com"
$ws1.Cells.Item(5, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(5, 6).Value = "delivery@example.com"
$ws1.Cells.Item(5, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(5, 8).Value = "2. 부자재비"
$ws1.Cells.Item(5, 9).Value = "2) 창호"
$ws1.Cells.Item(5, 10).Value = "기타"
$ws1.Cells.Item(5, 11).Value = "화기감시자 조끼(적색)"
$ws1.Cells.Item(5, 12).Value = "KS규격-4"
$ws1.Cells.Item(5, 13).Value = 2
$ws1.Cells.Item(5, 14).Value = 8000
$ws1.Cells.Item(5, 15).Value = 17600
$ws1.Cells.Item(5, 16).ClearContents()

# --- Row 6 ---
$ws1.Cells.Item(6, 1).Value = "2025-08-28"
$ws1.Cells.Item(6, 2).Value = "2025-10-16"
$ws1.Cells.Item(6, 3).Value = "티에스이앤씨"
$ws1.Cells.Item(6, 4).Value = "티에스이앤씨@example.com"
$ws1.Cells.Item(6, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(6, 6).Value = "delivery@example.com"
$ws1.Cells.Item(6, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(6, 8).Value = "2. 부자재비"
$ws1.Cells.Item(6, 9).Value = "2) 창호"
$ws1.Cells.Item(6, 10).Value = "기타"
$ws1.Cells.Item(6, 11).Value = "칼블럭 8*70"
$ws1.Cells.Item(6, 12).Value = "KS규격-5"
$ws1.Cells.Item(6, 13).Value = 1000
$ws1.Cells.Item(6, 14).Value = 119
$ws1.Cells.Item(6, 15).Value = 130900
$ws1.Cells.Item(6, 16).ClearContents()

# --- Row 7 ---
$ws1.Cells.Item(7, 1).Value = "2025-08-27"
$ws1.Cells.Item(7, 2).Value = "2025-10-08"
$ws1.Cells.Item(7, 3).Value = "티에스이앤씨"
$ws1.Cells.Item(7, 4).Value = "티에스이앤씨@example.com"
$ws1.Cells.Item(7, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(7, 6).Value = "delivery@example.com"
$ws1.Cells.Item(7, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(7, 8).Value = "2. 부자재비"
$ws1.Cells.Item(7, 9).Value = "2) 창호"
$ws1.Cells.Item(7, 10).Value = "기타"
$ws1.Cells.Item(7, 11).Value = "2차 - 타격공구"
$ws1.Cells.Item(7, 12).Value = "KS규격-6"
$ws1.Cells.Item(7, 13).Value = 2
$ws1.Cells.Item(7, 14).Value = 17500
$ws1.Cells.Item(7, 15).Value = 38500
$ws1.Cells.Item(7, 16).Value = "2차"

# Restore Normal style on A2:B7 now that the text values are safely stored
# (keeps the final cell style identical to the rest of the untouched data rows).
$ws1.Range("A2:B7").Style = "Normal"

# --- Remove old column Q entirely (also shrinks dimension down to A1:P7) ---
$ws1.Columns.Item(17).Delete()

# --- Sheets 갑지 (2) and 을지 (3): drop empty remark placeholder cells I2:I6 ---
foreach ($sheetIdx in 2, 3) {
  $ws = $wb.Worksheets.Item($sheetIdx)
  for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 9).ClearContents()
  }
}
